$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "before merge"
$ws.Range("E13").Value = "after"

$ws.Range("E13").Select()
